$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet ("Extra") so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ExposureGroup <-> Portname"

# Header row
$ws.Range("A1").Value = "ExposureGroup"
$ws.Range("B1").Value = "Portname"

# Data rows
$ws.Range("A2").Value = "USHU_Hurricane_Full"
$ws.Range("B2").Value = "USHU"

$ws.Range("A3").Value = "USEQ"
$ws.Range("B3").Value = "USEQ"

$ws.Range("A4").Value = "USFL_Flood_Full"
$ws.Range("B4").Value = "USFL"

$ws.Range("A5").Value = "CBHU_Hurricane_Full"
$ws.Range("B5").Value = "CBHU"

$ws.Range("A6").Value = "CBEQ_QuakeBC"
$ws.Range("B6").Value = "CBEQ"

$ws.Range("A7").Value = "MEHU"
$ws.Range("B7").Value = "MEHU"

$ws.Range("A8").Value = "PRHU"
$ws.Range("B8").Value = "PRHU"

# Style the header row to match the look used on the other sheets' header rows:
# bold text, centered horizontally, top-aligned vertically, thin box border.
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
